$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 "I0" and J1 "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), keyed by row number
$data = @{
    2 = @(7,7)
    3 = @(2,4)
    4 = @(9,9)
    5 = @(7,7)
    6 = @(8,8)
    7 = @(9,9)
    8 = @(8,8)
    9 = @(9,9)
    10 = @(7,7)
    11 = @(7,7)
    12 = @(9,9)
    13 = @(9,9)
    14 = @(7,7)
    15 = @(9,9)
    16 = @(7,8)
    17 = @(8,8)
    18 = @(9,9)
    19 = @(6,6)
    20 = @(8,8)
    21 = @(9,9)
    22 = @(8,8)
    23 = @(8,8)
    24 = @(6,7)
    25 = @(9,9)
    26 = @(7,8)
    27 = @(10,10)
    28 = @(8,8)
    29 = @(7,7)
    30 = @(5,5)
    31 = @(7,7)
    32 = @(8,8)
    33 = @(8,8)
    34 = @(7,7)
    35 = @(8,8)
    36 = @(7,7)
    37 = @(10,10)
    38 = @(9,9)
    39 = @(1,2)
    40 = @(7,7)
    41 = @(8,8)
    42 = @(9,9)
    43 = @(8,8)
    44 = @(7,7)
    45 = @(7,7)
    46 = @(7,7)
    47 = @(8,8)
    48 = @(7,7)
    49 = @(7,7)
    50 = @(7,7)
    51 = @(8,9)
    52 = @(6,6)
    53 = @(7,8)
    54 = @(8,8)
    55 = @(9,9)
    56 = @(8,9)
    57 = @(8,8)
    58 = @(9,9)
    59 = @(7,7)
    60 = @(9,9)
    61 = @(7,7)
    62 = @(8,9)
    63 = @(9,9)
    64 = @(9,9)
    65 = @(6,7)
    66 = @(7,7)
    67 = @(7,7)
    68 = @(6,6)
    69 = @(7,7)
    70 = @(8,8)
    71 = @(9,9)
    72 = @(7,7)
    73 = @(8,8)
    74 = @(6,6)
    75 = @(7,7)
    76 = @(6,6)
    77 = @(3,3)
    78 = @(3,3)
    79 = @(4,4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item([int]$row, 9).Value = $vals[0]
    $ws.Cells.Item([int]$row, 10).Value = $vals[1]
}

# Copy header formatting (style) from H1 onto the new I1:J1 headers,
# reusing the existing bold/border/center style rather than creating a new one
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

Write-Host "Added I0/IF columns"